# BIS-1002: removed "Internal Assignment" column from export.
# The "Internal Assignment" header (O4) and its per-row "FALSE" values
# (O5:O7) are cleared out, which also drops the now-unused
# "Internal Assignment" entry from the shared strings table on save.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Range("O4:O7").ClearContents()
